# Update the Metadata sheet: URL and Date values
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mindfulness-reminder-time"
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# Refresh the Elements sheet best-fit column widths (re-generated doc build narrowed them)
$elements = $wb.Worksheets.Item("Elements")
$elements.Columns.Item(1).ColumnWidth = 15.666666666666666
$elements.Columns.Item(2).ColumnWidth = 15.666666666666666
$elements.Columns.Item(3).ColumnWidth = 9.0
$elements.Columns.Item(3).Hidden = $true
$elements.Columns.Item(4).ColumnWidth = 6.166666666666667
$elements.Columns.Item(4).Hidden = $true
$elements.Columns.Item(5).ColumnWidth = 4.5
$elements.Columns.Item(6).ColumnWidth = 3.1666666666666665
$elements.Columns.Item(7).ColumnWidth = 3.5
$elements.Columns.Item(8).ColumnWidth = 11.833333333333334
$elements.Columns.Item(9).ColumnWidth = 9.666666666666666
$elements.Columns.Item(10).ColumnWidth = 19.833333333333332
$elements.Columns.Item(11).ColumnWidth = 7.5
$elements.Columns.Item(12).ColumnWidth = 99.83333333333333
$elements.Columns.Item(13).ColumnWidth = 99.83333333333333
$elements.Columns.Item(14).ColumnWidth = 99.83333333333333
$elements.Columns.Item(15).ColumnWidth = 11.5
$elements.Columns.Item(16).ColumnWidth = 19.833333333333332
$elements.Columns.Item(17).ColumnWidth = 19.833333333333332
$elements.Columns.Item(18).ColumnWidth = 19.833333333333332
$elements.Columns.Item(19).ColumnWidth = 19.833333333333332
$elements.Columns.Item(20).ColumnWidth = 7.0
$elements.Columns.Item(21).ColumnWidth = 12.833333333333334
$elements.Columns.Item(22).ColumnWidth = 13.166666666666666
$elements.Columns.Item(23).ColumnWidth = 14.166666666666666
$elements.Columns.Item(24).ColumnWidth = 13.833333333333334
$elements.Columns.Item(25).ColumnWidth = 16.166666666666668
$elements.Columns.Item(26).ColumnWidth = 14.333333333333334
$elements.Columns.Item(27).ColumnWidth = 4.166666666666667
$elements.Columns.Item(28).ColumnWidth = 17.166666666666668
$elements.Columns.Item(29).ColumnWidth = 33.666666666666664
$elements.Columns.Item(30).ColumnWidth = 12.666666666666666
$elements.Columns.Item(31).ColumnWidth = 10.5
$elements.Columns.Item(31).Hidden = $true
$elements.Columns.Item(32).ColumnWidth = 14.166666666666666
$elements.Columns.Item(32).Hidden = $true
$elements.Columns.Item(33).ColumnWidth = 7.333333333333333
$elements.Columns.Item(33).Hidden = $true
$elements.Columns.Item(34).ColumnWidth = 7.666666666666667
$elements.Columns.Item(35).ColumnWidth = 99.83333333333333
$elements.Columns.Item(37).ColumnWidth = 18.666666666666668
